$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.393.38"
$ws.Range("E2").Value = "  -4.56%  "
$ws.Range("D3").Value = "2.957.42"
$ws.Range("E3").Value = "  -6.38%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "538.54"
$ws.Range("E5").Value = "  -5.36%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "150.90"
$ws.Range("E6").Value = "  -7.52%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.999"
$ws.Range("E7").Value = "  -0.18%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.561"
$ws.Range("E8").Value = "  -3.85%  "
$ws.Range("D9").Value = "2.963.93"
$ws.Range("E9").Value = "  -6.11%  "
$ws.Range("E10").Value = "  -4.10%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.07"
$ws.Range("E11").Value = "  -8.53%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.364"
$ws.Range("E12").Value = "  -5.04%  "
$ws.Range("D13").Value = "3.466.97"
$ws.Range("E13").Value = "  -6.52%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.123"
$ws.Range("E14").Value = "  -3.13%  "
$ws.Range("D15").Value = "61.378.19"
$ws.Range("E15").Value = "  -4.72%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "23.57"
$ws.Range("E16").Value = "  -6.69%  "
$ws.Range("D17").Value = "2.956.64"
$ws.Range("E17").Value = "  -6.23%  "
$ws.Range("E18").Value = "  -6.60%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "5.11"
$ws.Range("E19").Value = "  -3.21%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "380.73"
$ws.Range("E20").Value = "  -6.31%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.90"
$ws.Range("E21").Value = "  -6.24%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.64"
$ws.Range("E22").Value = "  -6.23%  "
$ws.Range("E23").Value = "  -0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "64.89"
$ws.Range("E24").Value = "  -5.18%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.468"
$ws.Range("E25").Value = "  -3.01%  "
$ws.Range("D26").Value = "3.079.03"
$ws.Range("E26").Value = "  -6.74%  "
$ws.Range("E27").Value = "  -5.80%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "0.995"
$ws.Range("E28").Value = "  -0.49%  "
$ws.Range("D29").Value = "0.0₃0921"
$ws.Range("E29").Value = "  -10.33%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "8.24"
$ws.Range("E30").Value = "  -7.00%  "
$ws.Range("E31").Value = "  -0.03%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "1.71"
$ws.Range("E32").Value = "  -5.55%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "20.31"
$ws.Range("E33").Value = "  -4.31%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "158.68"
$ws.Range("E34").Value = "  +1.91%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "5.97"
$ws.Range("E35").Value = "  -5.71%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "4.60"
$ws.Range("E36").Value = "  -5.68%  "
$ws.Range("E37").Value = "  -6.68%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.27"
$ws.Range("E38").Value = "  -5.67%  "
$ws.Range("E39").Value = "  -9.10%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.90"
$ws.Range("E40").Value = "  -4.42%  "
$ws.Range("D41").Value = "2.406.07"
$ws.Range("E41").Value = "  -10.40%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "22.07"
$ws.Range("E42").Value = "  -8.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "36.78"
$ws.Range("E43").Value = "  -4.63%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.662"
$ws.Range("E44").Value = "  -4.90%  "
$ws.Range("E45").Value = "  -4.74%  "
$ws.Range("E46").Value = "  -0.27%  "
$ws.Range("E47").Value = "  -5.17%  "
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "4.86"
$ws.Range("E48").Value = "  -10.38%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0946"
$ws.Range("E49").Value = "  -3.84%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "10.47"
$ws.Range("E50").Value = "  -0.26%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "19.53"
$ws.Range("E51").Value = "  -9.11%  "
